# Daily attendance processing - 2026-01-01 14:58:56
# Reverse the order of comma-separated entries in the "Recorded By" (column G)
# values wherever more than one entry is present. Single-entry cells are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversedParts = $parts[($parts.Length - 1)..0]
        $newVal = [string]::Join(", ", $reversedParts)
        $cell.Value = $newVal
    }
}
